$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$headerRange = $ws.Range("A1:U1")

# Backup original formatting into an off-screen range first (copy format only)
$backupRange = $ws.Range("A100:U100")
$headerRange.Copy()
$backupRange.PasteSpecial(-4122) # xlPasteFormats
Write-Host "backed up formats"

$headerRange.ClearFormats()
$range = $ws.Range("A1:U61")
$tbl = $ws.ListObjects.Add(1, $range, $null, 1)
Write-Host "table created"

$backupRange.Copy()
$headerRange.PasteSpecial(-4122) # xlPasteFormats
Write-Host "restored formats from backup"

$backupRange.Clear()
Write-Host "cleared backup range"
